$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D10").Value = 1170.388941922528
$ws.Range("D11").Value = 1170.388941922528
$ws.Range("D12").Value = 1056.008723141073
$ws.Range("D13").Value = 1056.008723141073
